# Auto-generated edit script applying the Behemoth_Profits market-data refresh diff.
# Updates currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ
# columns (H, I, J, K, L, M, N) on specific rows across all 8 worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 602.61536
$ws.Range("I41").Value = 648.6667
$ws.Range("K41").Value = 648.6667
$ws.Range("M41").Value = -208.6667
$ws.Range("H64").Value = 4795.364
$ws.Range("J64").Value = 4999.8945
$ws.Range("L64").Value = 4999.8945
$ws.Range("N64").Value = -5495.8945
$ws.Range("H67").Value = 4795.364
$ws.Range("J67").Value = 4999.8945
$ws.Range("L67").Value = 4999.8945
$ws.Range("N67").Value = -6715.8945
$ws.Range("H93").Value = 124998.5
$ws.Range("J93").Value = 124998.5
$ws.Range("L93").Value = 124998.5
$ws.Range("N93").Value = -129990.5
$ws.Range("H126").Value = 89500
$ws.Range("J126").Value = 133000
$ws.Range("L126").Value = 133000
$ws.Range("N126").Value = -142880
$ws.Range("H127").Value = 7231.091
$ws.Range("I127").Value = 757
$ws.Range("J127").Value = 15000
$ws.Range("K127").Value = 2271
$ws.Range("L127").Value = 45000
$ws.Range("M127").Value = 2689
$ws.Range("N127").Value = -54920
$ws.Range("H132").Value = 604.97144
$ws.Range("I132").Value = 610.6418
$ws.Range("K132").Value = 1831.9254
$ws.Range("M132").Value = 698.0745999999999
$ws.Range("H137").Value = 375851.56
$ws.Range("I137").Value = 717019.5
$ws.Range("J137").Value = 8439.923000000001
$ws.Range("K137").Value = 2151058.5
$ws.Range("L137").Value = 25319.769
$ws.Range("M137").Value = -2148508.5
$ws.Range("N137").Value = -30419.769
$ws.Range("H138").Value = 2420.25
$ws.Range("J138").Value = 2525.0334
$ws.Range("L138").Value = 7575.100199999999
$ws.Range("N138").Value = -17855.1002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 45459000
$ws.Range("J45").Value = 6953.25
$ws.Range("L45").Value = 6953.25
$ws.Range("N45").Value = -7707.25
$ws.Range("H61").Value = 28850534
$ws.Range("J61").Value = 19235408
$ws.Range("L61").Value = 19235408
$ws.Range("N61").Value = -19235832
$ws.Range("H74").Value = 13170319
$ws.Range("I74").Value = 25003738
$ws.Range("K74").Value = 25003738
$ws.Range("M74").Value = -25002864
$ws.Range("H77").Value = 13170319
$ws.Range("I77").Value = 25003738
$ws.Range("K77").Value = 125018690
$ws.Range("M77").Value = -125014322
$ws.Range("H106").Value = 47767.6
$ws.Range("J106").Value = 47767.6
$ws.Range("L106").Value = 47767.6
$ws.Range("N106").Value = -50291.6
$ws.Range("H132").Value = 12351505
$ws.Range("I132").Value = 19611292
$ws.Range("K132").Value = 58833876
$ws.Range("M132").Value = -58831346
$ws.Range("H136").Value = 28850534
$ws.Range("J136").Value = 19235408
$ws.Range("L136").Value = 57706224
$ws.Range("N136").Value = -57711324

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2401.8096
$ws.Range("I94").Value = 3374.25
$ws.Range("J94").Value = 2173
$ws.Range("K94").Value = 3374.25
$ws.Range("L94").Value = 2173
$ws.Range("M94").Value = -2923.25
$ws.Range("N94").Value = -3075
$ws.Range("H106").Value = 63513.75
$ws.Range("J106").Value = 63513.75
$ws.Range("L106").Value = 63513.75
$ws.Range("N106").Value = -66037.75
$ws.Range("H134").Value = 3404368.5
$ws.Range("I134").Value = 1810.1923
$ws.Range("K134").Value = 5430.5769
$ws.Range("M134").Value = -2895.5769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1798436.8
$ws.Range("I31").Value = 2637.6
$ws.Range("J31").Value = 2920811.2
$ws.Range("K31").Value = 2637.6
$ws.Range("L31").Value = 2920811.2
$ws.Range("M31").Value = -2342.6
$ws.Range("N31").Value = -2921401.2
$ws.Range("H34").Value = 1798436.8
$ws.Range("I34").Value = 2637.6
$ws.Range("J34").Value = 2920811.2
$ws.Range("K34").Value = 2637.6
$ws.Range("L34").Value = 2920811.2
$ws.Range("M34").Value = -2435.6
$ws.Range("N34").Value = -2921215.2
$ws.Range("H58").Value = 13842.333
$ws.Range("I58").Value = 10756.5
$ws.Range("K58").Value = 10756.5
$ws.Range("M58").Value = -10553.5
$ws.Range("H122").Value = 2255
$ws.Range("I122").Value = 2303.2856
$ws.Range("K122").Value = 6909.8568
$ws.Range("M122").Value = -4459.8568
$ws.Range("H136").Value = 13842.333
$ws.Range("I136").Value = 10756.5
$ws.Range("K136").Value = 32269.5
$ws.Range("M136").Value = -29719.5
$ws.Range("H141").Value = 300822.94
$ws.Range("J141").Value = 313436.88
$ws.Range("L141").Value = 313436.88
$ws.Range("N141").Value = -323796.88

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 964274.5600000001
$ws.Range("I4").Value = 1333540.4
$ws.Range("K4").Value = 4000621.2
$ws.Range("M4").Value = -4000509.2
$ws.Range("H38").Value = 125.55556
$ws.Range("I38").Value = 149
$ws.Range("K38").Value = 447
$ws.Range("M38").Value = -100
$ws.Range("H56").Value = 6250.9
$ws.Range("I56").Value = 6250.9
$ws.Range("K56").Value = 6250.9
$ws.Range("M56").Value = -5720.9
$ws.Range("H80").Value = 4013.6365
$ws.Range("J80").Value = 4013.6365
$ws.Range("L80").Value = 12040.9095
$ws.Range("N80").Value = -13912.9095
$ws.Range("H83").Value = 4013.6365
$ws.Range("J83").Value = 4013.6365
$ws.Range("L83").Value = 36122.7285
$ws.Range("N83").Value = -45482.7285
$ws.Range("H139").Value = 4765.04
$ws.Range("J139").Value = 4991.6924
$ws.Range("L139").Value = 14975.0772
$ws.Range("N139").Value = -25255.0772
$ws.Range("H140").Value = 218672.72
$ws.Range("I140").Value = 218672.72
$ws.Range("K140").Value = 656018.16
$ws.Range("M140").Value = -650838.16

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 26382728
$ws.Range("J11").Value = 35002332
$ws.Range("L11").Value = 35002332
$ws.Range("N11").Value = -35002610
$ws.Range("H70").Value = 6323.6
$ws.Range("I70").Value = 5991.8
$ws.Range("J70").Value = 6655.4
$ws.Range("K70").Value = 5991.8
$ws.Range("L70").Value = 6655.4
$ws.Range("M70").Value = -5721.8
$ws.Range("N70").Value = -7195.4
$ws.Range("H73").Value = 6323.6
$ws.Range("I73").Value = 5991.8
$ws.Range("J73").Value = 6655.4
$ws.Range("K73").Value = 5991.8
$ws.Range("L73").Value = 6655.4
$ws.Range("M73").Value = -5055.8
$ws.Range("N73").Value = -8527.4
$ws.Range("H132").Value = 22729496
$ws.Range("I132").Value = 24392530
$ws.Range("J132").Value = 1366.3334
$ws.Range("K132").Value = 73177590
$ws.Range("L132").Value = 4099.0002
$ws.Range("M132").Value = -73175060
$ws.Range("N132").Value = -9159.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 58824020
$ws.Range("I55").Value = 66667156
$ws.Range("J55").Value = 499.5
$ws.Range("K55").Value = 66667156
$ws.Range("L55").Value = 499.5
$ws.Range("M55").Value = -66666983
$ws.Range("N55").Value = -845.5
$ws.Range("H68").Value = 4500
$ws.Range("J68").Value = 4500
$ws.Range("L68").Value = 4500
$ws.Range("N68").Value = -5998
$ws.Range("H71").Value = 4500
$ws.Range("J71").Value = 4500
$ws.Range("L71").Value = 22500
$ws.Range("N71").Value = -29988
$ws.Range("H100").Value = 2044
$ws.Range("I100").Value = 1430
$ws.Range("J100").Value = 4500
$ws.Range("K100").Value = 1430
$ws.Range("L100").Value = 4500
$ws.Range("M100").Value = -889
$ws.Range("N100").Value = -5582
$ws.Range("H132").Value = 2832.1702
$ws.Range("I132").Value = 2848.8125
$ws.Range("J132").Value = 2796.6667
$ws.Range("K132").Value = 8546.4375
$ws.Range("L132").Value = 8390.000100000001
$ws.Range("M132").Value = -6016.4375
$ws.Range("N132").Value = -13450.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H136").Value = 4416.6333
$ws.Range("J136").Value = 15447.5
$ws.Range("L136").Value = 46342.5
$ws.Range("N136").Value = -51442.5

